$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.6380968573760335
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 0.6112767687720249
$ws.Range("H2").Value = 0.4554673721340389
$ws.Range("D3").Value = 0.5732536407371513
$ws.Range("G3").Value = 0.6233937446395733
$ws.Range("H3").Value = 0.3334914611005693
$ws.Range("D4").Value = 0.6150516029863857
$ws.Range("G4").Value = 0.6099074018035465
$ws.Range("H4").Value = 0.4999999999999998
$ws.Range("D5").Value = 0.5808341423387191
$ws.Range("F5").Value = 42
$ws.Range("G5").Value = 0.5588326722920717
$ws.Range("H5").Value = 0.5833333333333334
$ws.Range("D6").Value = 0.5146467438650101
$ws.Range("G6").Value = 0.6099074018035465
$ws.Range("H6").Value = 0.5

$ws.Range("I2").Value = '[1 1 1 0 1 0 1 0 1 0 1 1 1 0 0 0 0 1 0 0 0 0 1 0]'
$ws.Range("J2").Value = '[0 0 0 1 0 1 1 1 0 0 0 1 1 1 1 0 0 1 0 0 1 0 0 0]'
$ws.Range("I3").Value = '[0 1 1 0 1 0 0 0 1 1 1 0 0 0 1 0 1 0 1 1 0 0 1 0]'
$ws.Range("J3").Value = '[0 0 0 0 0 0 1 0 0 0 0 1 0 1 0 1 0 0 1 0 0 0 0 1]'
$ws.Range("I4").Value = '[0 1 1 0 0 1 0 0 0 1 0 1 0 1 0 1 0 0 1 0 0 1 1 1]'
$ws.Range("J4").Value = '[1 0 0 1 1 1 0 0 0 0 1 1 0 1 1 1 0 0 0 1 1 1 0 1]'
$ws.Range("I5").Value = '[1 0 1 0 0 0 1 1 1 1 1 1 1 0 0 0 0 0 0 0 1 1 0 0]'
$ws.Range("J5").Value = '[1 0 0 0 0 1 0 1 1 1 0 0 1 1 0 0 0 1 0 1 0 1 0 1]'
$ws.Range("I6").Value = '[1 0 1 0 1 1 0 0 0 1 0 1 0 1 1 1 0 1 0 1 0 0 0 0]'
$ws.Range("J6").Value = '[0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 1 0 1 0 1 1 0 1 1]'

$B2_val = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('MinMaxScaler',
                                                  MinMaxScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta...
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f5b614fa6d0>),
                ('model',
                 AdaBoostClassifier(estimator=LogisticRegression(class_weight='balanced',
                                                                 l1_ratio=0.01,
                                                                 max_iter=1000,
                                                                 penalty='elasticnet',
                                                                 random_state=42,
                                                                 solver='saga'),
                                    n_estimators=10, random_state=42))])
'@
$ws.Range("B2").Value = $B2_val

$C2_val = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f5b613577c0>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('MinMaxScaler', MinMaxScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 10, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'elasticnet', 'model__estimator__l1_ratio': 0.01, 'model__estimator__class_weight': 'balanced'}
'@
$ws.Range("C2").Value = $C2_val

$B3_val = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('MinMaxScaler',
                                                  MinMaxScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', None),
                ('model',
                 AdaBoostClassifier(estimator=LogisticRegression(class_weight='balanced',
                                                                 l1_ratio=0.001,
                                                                 max_iter=1000,
                                                                 penalty='elasticnet',
                                                                 random_state=42,
                                                                 solver='saga'),
                                    n_estimators=10, random_state=42))])
'@
$ws.Range("B3").Value = $B3_val

$C3_val = @'
{'selector': None, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('MinMaxScaler', MinMaxScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 10, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'elasticnet', 'model__estimator__l1_ratio': 0.001, 'model__estimator__class_weight': 'balanced'}
'@
$ws.Range("C3").Value = $C3_val

$B4_val = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('MinMaxScaler',
                                                  MinMaxScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta...
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f5b61357a90>),
                ('model',
                 AdaBoostClassifier(estimator=LogisticRegression(class_weight='balanced',
                                                                 l1_ratio=0.001,
                                                                 max_iter=1000,
                                                                 penalty='elasticnet',
                                                                 random_state=42,
                                                                 solver='saga'),
                                    n_estimators=5, random_state=42))])
'@
$ws.Range("B4").Value = $B4_val

$C4_val = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f5b6163a0a0>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('MinMaxScaler', MinMaxScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 5, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'elasticnet', 'model__estimator__l1_ratio': 0.001, 'model__estimator__class_weight': 'balanced'}
'@
$ws.Range("C4").Value = $C4_val

$B5_val = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('MinMaxScaler',
                                                  MinMaxScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta...
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f5b6163a1f0>),
                ('model',
                 AdaBoostClassifier(estimator=LogisticRegression(class_weight='balanced',
                                                                 l1_ratio=0.01,
                                                                 max_iter=1000,
                                                                 penalty='elasticnet',
                                                                 random_state=42,
                                                                 solver='saga'),
                                    n_estimators=10, random_state=42))])
'@
$ws.Range("B5").Value = $B5_val

$C5_val = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f5ae82edc10>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('MinMaxScaler', MinMaxScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 10, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'elasticnet', 'model__estimator__l1_ratio': 0.01, 'model__estimator__class_weight': 'balanced'}
'@
$ws.Range("C5").Value = $C5_val

$B6_val = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('MinMaxScaler',
                                                  MinMaxScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta...
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f5b616be820>),
                ('model',
                 AdaBoostClassifier(estimator=LogisticRegression(class_weight='balanced',
                                                                 l1_ratio=0.01,
                                                                 max_iter=1000,
                                                                 penalty='elasticnet',
                                                                 random_state=42,
                                                                 solver='saga'),
                                    n_estimators=100, random_state=42))])
'@
$ws.Range("B6").Value = $B6_val

$C6_val = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f5b619783a0>, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('MinMaxScaler', MinMaxScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__n_estimators': 100, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'elasticnet', 'model__estimator__l1_ratio': 0.01, 'model__estimator__class_weight': 'balanced'}
'@
$ws.Range("C6").Value = $C6_val

